$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Breaking News, World News and Video from Al Jazeera"
$ws.Range("C2").Value = "News, analysis from the Middle East & worldwide, multimedia & interactives, opinions, documentaries, podcasts, long reads and broadcast"
$ws.Range("E2").Value = 0

# Row 3
$ws.Range("A3").Value = "Israel's war on Gaza live news: Attacks on besieged enclave kill 62 ..."
$ws.Range("C3").Value = "Israel, a major recipient of US military assistance for decades, is still due to receive billions of dollars of US aid and weaponry. “The"

# Row 4
$ws.Range("A4").Value = "Economy | Today's latest from Al Jazeera"
$ws.Range("B4").Value = 45433
$ws.Range("C4").Value = "Russian court seizes two European banks' assets amid Western sanctions. Freezing hundreds of billions of dollars in lenders' assets was part of dispute over gas"
$ws.Range("E4").Value = 1

# Row 5
$ws.Range("A5").Value = "Israel's war on Gaza live news: Deadly combat rages as Rafah ..."
$ws.Range("C5").Value = "Israel, a major recipient of US military assistance for decades, is still due to receive billions of dollars of US aid and weaponry. “The"

# Row 6
$ws.Range("A6").Value = "Gaza war: What does victory look like for the US and Israel? | Israel ..."
$ws.Range("B6").Value = 45433
$ws.Range("C6").Value = "Israel has said it is seeking an “absolute victory” over Hamas, as it continues to receive billions of dollars in unconditional military aid"

# Row 13
$ws.Range("B13").Value = 45431
